$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1463414634146341
$ws.Range("C2").Value = 0.6376306620209059
$ws.Range("J2").Value = 0.006968641114982578
$ws.Range("P2").Value = 0.132404181184669
$ws.Range("S2").Value = 0.07665505226480836
$ws.Range("C3").Value = 0.07035175879396985
$ws.Range("J3").Value = 0.03517587939698492
$ws.Range("P3").Value = 0.7587939698492462
$ws.Range("S3").Value = 0.135678391959799
$ws.Range("J4").Value = 0.06896551724137931
$ws.Range("P4").Value = 0.7241379310344828
$ws.Range("S4").Value = 0.2068965517241379
$ws.Range("B6").Value = 0.04629629629629629
$ws.Range("D6").Value = 0.01388888888888889
$ws.Range("F6").Value = 0.04629629629629629
$ws.Range("J6").Value = 0.2175925925925926
$ws.Range("O6").Value = 0.009259259259259259
$ws.Range("Q6").Value = 0.1527777777777778
$ws.Range("R6").Value = 0.1111111111111111
$ws.Range("S6").Value = 0.4027777777777778
$ws.Range("B7").Value = 0.1565656565656566
$ws.Range("D7").Value = 0.01515151515151515
$ws.Range("F7").Value = 0.03535353535353535
$ws.Range("J7").Value = 0.1060606060606061
$ws.Range("O7").Value = 0.02525252525252525
$ws.Range("Q7").Value = 0.1515151515151515
$ws.Range("R7").Value = 0.1060606060606061
$ws.Range("S7").Value = 0.404040404040404
$ws.Range("B8").Value = 0.1133501259445844
$ws.Range("D8").Value = 0.01007556675062972
$ws.Range("F8").Value = 0.08816120906801007
$ws.Range("J8").Value = 0.1007556675062972
$ws.Range("O8").Value = 0.01763224181360202
$ws.Range("Q8").Value = 0.1536523929471033
$ws.Range("R8").Value = 0.1057934508816121
$ws.Range("S8").Value = 0.4105793450881612
$ws.Range("B9").Value = 0.0718562874251497
$ws.Range("D9").Value = 0.01197604790419162
$ws.Range("F9").Value = 0.09580838323353294
$ws.Range("J9").Value = 0.1017964071856287
$ws.Range("O9").Value = 0.005988023952095809
$ws.Range("Q9").Value = 0.1497005988023952
$ws.Range("R9").Value = 0.155688622754491
$ws.Range("S9").Value = 0.407185628742515
$ws.Range("B10").Value = 0.1182707993474715
$ws.Range("D10").Value = 0.01468189233278956
$ws.Range("E10").Value = 0.002446982055464926
$ws.Range("F10").Value = 0.0864600326264274
$ws.Range("J10").Value = 0.1060358890701468
$ws.Range("O10").Value = 0.009787928221859706
$ws.Range("Q10").Value = 0.1924959216965742
$ws.Range("R10").Value = 0.09216965742251224
$ws.Range("S10").Value = 0.3776508972267537
$ws.Range("G11").Value = 0.1423728813559322
$ws.Range("J11").Value = 0.07457627118644068
$ws.Range("K11").Value = 0.1864406779661017
$ws.Range("L11").Value = 0.5898305084745763
$ws.Range("S11").Value = 0.006779661016949152
$ws.Range("G12").Value = 0.768361581920904
$ws.Range("J12").Value = 0.1694915254237288
$ws.Range("K12").Value = 0.01694915254237288
$ws.Range("L12").Value = 0.02824858757062147
$ws.Range("S12").Value = 0.01694915254237288
$ws.Range("G13").Value = 0.5365853658536586
$ws.Range("J13").Value = 0.3658536585365854
$ws.Range("S13").Value = 0.0975609756097561
$ws.Range("F15").Value = 0.005494505494505495
$ws.Range("H15").Value = 0.1923076923076923
$ws.Range("I15").Value = 0.04395604395604396
$ws.Range("J15").Value = 0.3516483516483517
$ws.Range("K15").Value = 0.08791208791208792
$ws.Range("M15").Value = 0.005494505494505495
$ws.Range("O15").Value = 0.06593406593406594
$ws.Range("S15").Value = 0.2472527472527473
$ws.Range("F16").Value = 0.01477832512315271
$ws.Range("H16").Value = 0.1280788177339902
$ws.Range("I16").Value = 0.06896551724137931
$ws.Range("J16").Value = 0.4088669950738916
$ws.Range("M16").Value = 0.03940886699507389
$ws.Range("O16").Value = 0.0541871921182266
$ws.Range("S16").Value = 0.1428571428571428
$ws.Range("F17").Value = 0.01036269430051814
$ws.Range("H17").Value = 0.1709844559585492
$ws.Range("I17").Value = 0.09844559585492228
$ws.Range("J17").Value = 0.4404145077720207
$ws.Range("K17").Value = 0.1295336787564767
$ws.Range("M17").Value = 0.01813471502590673
$ws.Range("N17").Value = 0.002590673575129534
$ws.Range("O17").Value = 0.04922279792746114
$ws.Range("S17").Value = 0.08031088082901554
$ws.Range("H18").Value = 0.2053571428571428
$ws.Range("I18").Value = 0.08928571428571429
$ws.Range("J18").Value = 0.4196428571428572
$ws.Range("K18").Value = 0.08928571428571429
$ws.Range("M18").Value = 0.008928571428571428
$ws.Range("O18").Value = 0.05357142857142857
$ws.Range("S18").Value = 0.1339285714285714
$ws.Range("F19").Value = 0.01012658227848101
$ws.Range("H19").Value = 0.1915611814345992
$ws.Range("I19").Value = 0.07257383966244725
$ws.Range("J19").Value = 0.4210970464135021
$ws.Range("K19").Value = 0.1021097046413502
$ws.Range("M19").Value = 0.01940928270042194
$ws.Range("O19").Value = 0.06160337552742616
$ws.Range("S19").Value = 0.1215189873417721
